# Updated symbol list on Tue Dec 27 08:48:16 UTC 2022 with GitHub Actions
#
# Applies the per-coin price / label refresh captured in the diff:
#  - Column D ("Price") cells are stored as text in the source sheet, so we
#    force the Text number format before writing the new value. That keeps
#    things like trailing zeros ("0.0005940") and near-zero magnitudes
#    ("0.00008001") intact instead of Excel collapsing them to a float.
#  - Column B/C/E are plain text already (coin name / link / volume label)
#    so a normal .Value assignment is sufficient.
#  - Rows 41 and 43 had their BKEXToken/KickToken rows swapped (plus fresh
#    prices/labels), so each column on those two rows is rewritten in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - BNB
Set-TextValue $ws.Range("D2") "242.88"

# Row 3 - OKB
Set-TextValue $ws.Range("D3") "23.01"

# Row 4 - HuobiToken
Set-TextValue $ws.Range("D4") "5.404"

# Row 5 - Cronos
Set-TextValue $ws.Range("D5") "0.05958"

# Row 7
Set-TextValue $ws.Range("D7") "6.503"

# Row 8
Set-TextValue $ws.Range("D8") "0.8142"

# Row 9
Set-TextValue $ws.Range("D9") "0.9233"

# Row 10
Set-TextValue $ws.Range("D10") "0.1435"

# Row 11
Set-TextValue $ws.Range("D11") "0.07380"

# Row 12
Set-TextValue $ws.Range("D12") "0.03301"

# Row 13
Set-TextValue $ws.Range("D13") "0.03081"

# Row 14
Set-TextValue $ws.Range("D14") "0.09346"

# Row 15
Set-TextValue $ws.Range("D15") "3.845"

# Row 16
Set-TextValue $ws.Range("D16") "0.001584"

# Row 18 - One (ONE) : price + label (was Worstin24h -> now restored, adds suffix)
Set-TextValue $ws.Range("D18") "0.0005940"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19
Set-TextValue $ws.Range("D19") "0.005900"

# Row 20 - BitKan (KAN): price + label (drops Bestin24h suffix)
Set-TextValue $ws.Range("D20") "0.001262"
$ws.Range("E20").Value = "19BitKanKAN"

# Row 21
Set-TextValue $ws.Range("D21") "0.004846"

# Row 22
Set-TextValue $ws.Range("D22") "0.00008001"

# Row 23
Set-TextValue $ws.Range("D23") "3.568"

# Row 27
Set-TextValue $ws.Range("D27") "0.0002339"

# Row 40 - IDEX
Set-TextValue $ws.Range("D40") "0.03946"

# Row 41 - was BKEXToken, now KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006340"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 - CEJI : price + label (adds Bestin24h suffix)
Set-TextValue $ws.Range("D42") "0.003900"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"

# Row 43 - was KickToken, now BKEXToken
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D43") "0.1073"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# Row 44
Set-TextValue $ws.Range("D44") "0.008910"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005179"

# Row 47
Set-TextValue $ws.Range("D47") "0.6800"

# Row 48
Set-TextValue $ws.Range("D48") "0.002148"
